$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RMA Details Maintenance Grid")

# Row 3 corresponds to the "Pro-Stock Product (Stock Product No Track)" line.
# Update the Sales Order Line / Shipper Line / Id values to the new RMA test case
# (RMA-JBWV) generated for this flow (SO To inspection order SO to RMA Receipt).
$ws.Range("E3").Value = "RMA-JBWV-001"
$ws.Range("F3").Value = "RMA-JBWV-1-1"
$ws.Range("J3").Value = "a6h1K000000Q2JVQA0"
